$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.073.01'
$ws.Range('E2').Value = '  -0.59%  '
# Row 3
$ws.Range('D3').Value = '3.328.05'
$ws.Range('E3').Value = '  +0.27%  '
# Row 4
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.36%  '
# Row 5
$ws.Range('D5').Value = '582.15'
$ws.Range('E5').Value = '  +3.26%  '
# Row 6
$ws.Range('D6').Value = '185.24'
$ws.Range('E6').Value = '  -2.98%  '
# Row 7
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.02%  '
# Row 8
$ws.Range('D8').Value = '3.325.41'
$ws.Range('E8').Value = '  +0.44%  '
# Row 9
$ws.Range('D9').Value = '0.576'
$ws.Range('E9').Value = '  -2.57%  '
# Row 10
$ws.Range('D10').Value = '0.180'
$ws.Range('E10').Value = '  -3.90%  '
# Row 11
$ws.Range('D11').Value = '0.577'
$ws.Range('E11').Value = '  -2.23%  '
# Row 12
$ws.Range('D12').Value = '47.02'
$ws.Range('E12').Value = '  -2.23%  '
# Row 13
$ws.Range('D13').Value = '0.0000267'
$ws.Range('E13').Value = '  -2.05%  '
# Row 14
$ws.Range('D14').Value = '636.11'
$ws.Range('E14').Value = '  +3.24%  '
# Row 15
$ws.Range('D15').Value = '3.848.03'
# Row 16
$ws.Range('D16').Value = '8.50'
$ws.Range('E16').Value = '  -2.58%  '
# Row 17
$ws.Range('D17').Value = '66.101.09'
$ws.Range('E17').Value = '  -0.52%  '
# Row 18
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '17.94'
$ws.Range('E18').Value = '  -1.21%  '
# Row 19
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.117'
$ws.Range('E19').Value = '  -0.07%  '
# Row 20
$ws.Range('D20').Value = '3.310.35'
$ws.Range('E20').Value = '  +0.09%  '
# Row 21
$ws.Range('D21').Value = '11.06'
$ws.Range('E21').Value = '  -0.69%  '
# Row 22
$ws.Range('D22').Value = '0.897'
$ws.Range('E22').Value = '  -1.90%  '
# Row 23
$ws.Range('D23').Value = '17.82'
$ws.Range('E23').Value = '  -3.52%  '
# Row 24
$ws.Range('D24').Value = '5.04'
$ws.Range('E24').Value = '  -2.41%  '
# Row 25
$ws.Range('D25').Value = '100.59'
$ws.Range('E25').Value = '  -1.56%  '
# Row 26
$ws.Range('D26').Value = '3.99'
$ws.Range('E26').Value = '  -0.87%  '
# Row 27
$ws.Range('D27').Value = '2.75'
$ws.Range('E27').Value = '  -0.64%  '
# Row 28
$ws.Range('D28').Value = '9.49'
$ws.Range('E28').Value = '  -3.30%  '
# Row 29
$ws.Range('D29').Value = '31.12'
$ws.Range('E29').Value = '  +2.35%  '
# Row 30
$ws.Range('D30').Value = '8.45'
$ws.Range('E30').Value = '  -2.50%  '
# Row 31
$ws.Range('D31').Value = '6.70'
$ws.Range('E31').Value = '  -0.93%  '
# Row 32
$ws.Range('D32').Value = '595.07'
$ws.Range('E32').Value = '  +3.62%  '
# Row 33
$ws.Range('D33').Value = '3.84'
$ws.Range('E33').Value = '  -7.70%  '
# Row 34
$ws.Range('D34').Value = '10.98'
$ws.Range('E34').Value = '  -1.52%  '
# Row 35
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.105'
$ws.Range('E35').Value = '  -0.66%  '
# Row 36
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '3.844.30'
$ws.Range('E36').Value = '  +2.29%  '
# Row 37
$ws.Range('E37').Value = '  +0.13%  '
# Row 38
$ws.Range('D38').Value = '55.97'
$ws.Range('E38').Value = '  -2.63%  '
# Row 39
$ws.Range('D39').Value = '0.0₃0702'
$ws.Range('E39').Value = '  -4.77%  '
# Row 40
$ws.Range('E40').Value = '  -3.41%  '
# Row 41
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = '32.69'
$ws.Range('E41').Value = '  -5.12%  '
# Row 42
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '2.65'
$ws.Range('E42').Value = '  -4.16%  '
# Row 43
$ws.Range('D43').Value = '3.17'
$ws.Range('E43').Value = '  -5.47%  '
# Row 44
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').Value = '3.42'
$ws.Range('E44').Value = '  +4.53%  '
# Row 45
$ws.Range('D45').Value = '0.336'
$ws.Range('E45').Value = '  -2.14%  '
# Row 46
$ws.Range('D46').Value = '0.0412'
$ws.Range('E46').Value = '  -4.00%  '
# Row 47
$ws.Range('D47').Value = '3.07'
$ws.Range('E47').Value = '  -14.29%  '
# Row 48
$ws.Range('D48').Value = '0.128'
$ws.Range('E48').Value = '  -1.64%  '
# Row 49
$ws.Range('E49').Value = '  +0.28%  '
# Row 50
$ws.Range('D50').Value = '2.54'
$ws.Range('E50').Value = '  -2.66%  '
# Row 51
$ws.Range('D51').Value = '130.73'
$ws.Range('E51').Value = '  +6.90%  '
